$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.966.93"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.285.92"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.59%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.406"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "3.856.62"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").Value = "66.062.92"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "3.270.24"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "429.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "3.436.07"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.196"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").Value = "2.787.34"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.770"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0658"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.42%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "314.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.08%  "
